# Auto-generated edit script: updates market-price-derived profit columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
# Row 28
$ws.Range("H28").Value = 1267.2
$ws.Range("I28").Value = 1446.5
$ws.Range("J28").Value = 550
$ws.Range("K28").Value = 1446.5
$ws.Range("L28").Value = 550
$ws.Range("M28").Value = -961.5
$ws.Range("N28").Value = -1520
# Row 40
$ws.Range("H40").Value = 4643.4
$ws.Range("I40").Value = 1483.5
$ws.Range("J40").Value = 6750
$ws.Range("K40").Value = 1483.5
$ws.Range("L40").Value = 6750
$ws.Range("M40").Value = -1308.5
$ws.Range("N40").Value = -7100
# Row 132
$ws.Range("H132").Value = 2963.6191
$ws.Range("I132").Value = 1108.125
$ws.Range("K132").Value = 3324.375
$ws.Range("M132").Value = -794.375
# Row 138
$ws.Range("H138").Value = 3203.0715
$ws.Range("I138").Value = 815.2857
$ws.Range("K138").Value = 2445.8571
$ws.Range("M138").Value = 2694.1429

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3425.5217
$ws.Range("J32").Value = 499.5
$ws.Range("L32").Value = 499.5
$ws.Range("N32").Value = -1073.5
# Row 45
$ws.Range("H45").Value = 2608.9
$ws.Range("J45").Value = 3365
$ws.Range("L45").Value = 3365
$ws.Range("N45").Value = -4119
# Row 63
$ws.Range("H63").Value = 2750
$ws.Range("I63").Value = 2750
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2750
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2064
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 2750
$ws.Range("I66").Value = 2750
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 13750
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -10318
$ws.Range("N66").ClearContents()
# Row 97
$ws.Range("H97").Value = 937.4545000000001
$ws.Range("I97").Value = 620.875
$ws.Range("K97").Value = 620.875
$ws.Range("M97").Value = -124.875

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 8
$ws.Range("H8").Value = 687.25
$ws.Range("I8").Value = 687.25
$ws.Range("K8").Value = 687.25
$ws.Range("M8").Value = -547.25
# Row 82
$ws.Range("H82").Value = 4500
$ws.Range("I82").Value = 4500
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 4500
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -4117
$ws.Range("N82").ClearContents()
# Row 85
$ws.Range("H85").Value = 4500
$ws.Range("I85").Value = 4500
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 4500
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -3174
$ws.Range("N85").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 10
$ws.Range("H10").Value = 308.66666
$ws.Range("I10").Value = 308.66666
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 308.66666
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -169.66666
$ws.Range("N10").ClearContents()
# Row 31
$ws.Range("H31").Value = 10219
$ws.Range("I31").Value = 10219
$ws.Range("K31").Value = 10219
$ws.Range("M31").Value = -9924
# Row 34
$ws.Range("H34").Value = 10219
$ws.Range("I34").Value = 10219
$ws.Range("K34").Value = 10219
$ws.Range("M34").Value = -10017
# Row 50
$ws.Range("H50").Value = 25333.334
$ws.Range("I50").Value = 7000
$ws.Range("J50").Value = 29000
$ws.Range("K50").Value = 7000
$ws.Range("L50").Value = 29000
$ws.Range("M50").Value = -6375
$ws.Range("N50").Value = -30250
# Row 132
$ws.Range("H132").Value = 4687.647
$ws.Range("I132").Value = 4537
$ws.Range("K132").Value = 13611
$ws.Range("M132").Value = -11081

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 103
$ws.Range("H103").Value = 4512
$ws.Range("I103").Value = 4025
$ws.Range("J103").Value = 4999
$ws.Range("K103").Value = 12075
$ws.Range("L103").Value = 14997
$ws.Range("M103").Value = -11196
$ws.Range("N103").Value = -16755
# Row 114
$ws.Range("H114").Value = 1774.3334
$ws.Range("J114").Value = 2149.5
$ws.Range("L114").Value = 6448.5
$ws.Range("N114").Value = -12956.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 10500033
$ws.Range("I11").Value = 10500033
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 10500033
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -10499894
$ws.Range("N11").ClearContents()
# Row 122
$ws.Range("H122").Value = 7167.4707
$ws.Range("I122").Value = 7365.4375
$ws.Range("K122").Value = 22096.3125
$ws.Range("M122").Value = -19646.3125

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5300.4
$ws.Range("I7").Value = 5000.5
$ws.Range("J7").Value = 6500
$ws.Range("K7").Value = 5000.5
$ws.Range("L7").Value = 6500
$ws.Range("M7").Value = -4888.5
$ws.Range("N7").Value = -6724
# Row 29
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
# Row 40
$ws.Range("H40").Value = 3113.8572
$ws.Range("I40").Value = 3113.8572
$ws.Range("K40").Value = 3113.8572
$ws.Range("M40").Value = -2977.8572
# Row 46
$ws.Range("H46").Value = 957.8
$ws.Range("I46").Value = 863.3333
$ws.Range("K46").Value = 863.3333
$ws.Range("M46").Value = -675.3333
# Row 126
$ws.Range("H126").Value = 5300.4
$ws.Range("I126").Value = 5000.5
$ws.Range("J126").Value = 6500
$ws.Range("K126").Value = 15001.5
$ws.Range("L126").Value = 19500
$ws.Range("M126").Value = -12531.5
$ws.Range("N126").Value = -24440
# Row 132
$ws.Range("H132").Value = 8071.625
$ws.Range("I132").Value = 8071.625
$ws.Range("K132").Value = 24214.875
$ws.Range("M132").Value = -21684.875

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 3866.6667
$ws.Range("I81").Value = 3866.6667
$ws.Range("K81").Value = 7733.3334
$ws.Range("M81").Value = -6672.3334
# Row 84
$ws.Range("H84").Value = 3866.6667
$ws.Range("I84").Value = 3866.6667
$ws.Range("K84").Value = 38666.667
$ws.Range("M84").Value = -33362.667
# Row 107
$ws.Range("H107").Value = 1245.7142
$ws.Range("I107").Value = 544
$ws.Range("K107").Value = 1632
$ws.Range("M107").Value = 288
# Row 122
$ws.Range("H122").Value = 287070.72
$ws.Range("I122").Value = 287070.72
$ws.Range("K122").Value = 861212.1599999999
$ws.Range("M122").Value = -858762.1599999999
# Row 126
$ws.Range("H126").Value = 3994
$ws.Range("I126").Value = 3994
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 11982
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -9512
$ws.Range("N126").ClearContents()
